$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Update the "License" paragraph (bold title run + descriptive runs)
# ---------------------------------------------------------------------------
# Locate the paragraph that still carries the old license description text.
$licenseIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Contains("is based on")) {
        $licenseIndex = $i
        break
    }
}

# Replace the bold title run's text; scope the Find to the paragraph so it
# only touches this occurrence (the same text also appears elsewhere).
$licensePara = $d.Paragraphs.Item($licenseIndex)
$pr = $licensePara.Range
$pr.Find.ClearFormatting()
$pr.Find.Execute("關鍵詞 (Biblica)", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Biblica Study Notes (Key Terms)", 2) | Out-Null

# Re-acquire the paragraph (content length changed) and find the end of the
# bold title run so we know where the non-bold remainder begins.
$licensePara = $d.Paragraphs.Item($licenseIndex)
$pr2 = $licensePara.Range
$pr2.Find.ClearFormatting()
$pr2.Find.Execute("Biblica Study Notes (Key Terms)", $false, $false, $false, $false, $false, $true, 1, $false, `
    "", 0) | Out-Null
$afterTitle = $pr2.End

# Delete everything after the title up to (not including) the paragraph mark;
# this removes the old descriptive text plus both hyperlink fields.
$paraEnd = $licensePara.Range.End
$deleteRange = $d.Range($afterTitle, $paraEnd - 1)
$deleteRange.Text = ""

# Insert the new descriptive text as plain (non-bold) runs.
$part1 = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. "
$part2 = "Biblica Study Notes"
$part3 = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."

$ins1 = $d.Range($afterTitle, $afterTitle)
$ins1.Text = $part1
$fix1 = $d.Range($afterTitle, $ins1.End)
$fix1.Font.Bold = 0

$pos2 = $ins1.End
$ins2 = $d.Range($pos2, $pos2)
$ins2.Text = $part2
$fix2 = $d.Range($pos2, $ins2.End)
$fix2.Font.Bold = 0

$pos3 = $ins2.End
$ins3 = $d.Range($pos3, $pos3)
$ins3.Text = $part3
$fix3 = $d.Range($pos3, $ins3.End)
$fix3.Font.Bold = 0

# Re-create the trailing empty run that originally closed the paragraph, so
# paragraph structure (trailing empty run) is preserved.
$endPos = $ins3.End
$tailA = $d.Range($endPos, $endPos)
$tailA.Font.Bold = -1
$tailB = $d.Range($endPos, $endPos)
$tailB.Font.Bold = 0

# ---------------------------------------------------------------------------
# Edit 2: Remove the whole "License Information" heading paragraph
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.Trim()
    if ($t -eq "License Information") {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# Edit 3: Remove the whole "This PDF version is provided under the same
# license." paragraph
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.Trim()
    if ($t -eq "This PDF version is provided under the same license.") {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# Edit 4: Remove the whole paragraph holding the italic key-term list that
# starts with "神, 神的兒子, 神的兒子們"
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Contains([char]0x795E) -and $t.Contains(",") -and $t.Length -gt 10) {
        $p.Range.Delete()
        break
    }
}
